# "Generate Report for Handback" -- fills in the Latest Target File / Latest
# Handback File / Latest Handback DateTime columns (and the Status column) on
# the zh-cn and de-de report sheets now that both language handbacks are in,
# and reflects that on the Overview sheet too.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fda62e73e1f2f0abb6c8d2e3ad656dde4c82fbb5/e2e/"
$mdFile1 = "3d820aaf-3ea1-4374-b429-de7c5d8ce9b8.md"
$mdFile2 = "8c7df7c7-1f23-4760-b8c4-1e2fd9b39bc9.md"

# ---------------------------------------------------------------------
# Overview sheet: both language-status columns move to "Handed back"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusHandedBack
$zh.Range("C3").Value = $statusHandedBack

$zh.Range("I2").Value = $mdFile1
$zh.Range("J2").Value = "3d820aaf-3ea1-4374-b429-de7c5d8ce9b8.ae5d4d9b2001a9660ea9b2f6bcfafbb4e9940eb9.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-07 17:12:50"

$zh.Range("I3").Value = $mdFile2
$zh.Range("J3").Value = "8c7df7c7-1f23-4760-b8c4-1e2fd9b39bc9.3ce779e4edda2ae5f3bf558c5e60296128ef04e6.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-07 17:12:50"

$zh.Hyperlinks.Add($zh.Range("I2"), ($repoBase + $mdFile1), "", "", $mdFile1)
$zh.Hyperlinks.Add($zh.Range("I3"), ($repoBase + $mdFile2), "", "", $mdFile2)

$zh.Columns.Item(3).ColumnWidth = 29.16
$zh.Columns.Item(9).ColumnWidth = 39.16
$zh.Columns.Item(10).ColumnWidth = 39.16

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusHandedBack
$de.Range("C3").Value = $statusHandedBack

$de.Range("I2").Value = $mdFile1
$de.Range("J2").Value = "3d820aaf-3ea1-4374-b429-de7c5d8ce9b8.ae5d4d9b2001a9660ea9b2f6bcfafbb4e9940eb9.de-de.xlf"
$de.Range("K2").Value = "2016-09-07 17:12:59"

$de.Range("I3").Value = $mdFile2
$de.Range("J3").Value = "8c7df7c7-1f23-4760-b8c4-1e2fd9b39bc9.3ce779e4edda2ae5f3bf558c5e60296128ef04e6.de-de.xlf"
$de.Range("K3").Value = "2016-09-07 17:12:59"

$de.Hyperlinks.Add($de.Range("I2"), ($repoBase + $mdFile1), "", "", $mdFile1)
$de.Hyperlinks.Add($de.Range("I3"), ($repoBase + $mdFile2), "", "", $mdFile2)

$de.Columns.Item(3).ColumnWidth = 29.16
$de.Columns.Item(9).ColumnWidth = 39.16
$de.Columns.Item(10).ColumnWidth = 39.16

# ---------------------------------------------------------------------
# Overview sheet column widths follow the Status text getting longer
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.16
$overview.Columns.Item(6).ColumnWidth = 29.16
